$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

$ws.Range("G2").Value = "wait(3);`nvalidate1;`nSwitchApp(NATIVE_APP);`nClickNativeIcon(VT200_0851_home_xpath);`nSwitchApp(WEBVIEW);`nlink_Click(signal_test_link);`nvalidate2;`nSelectTestToRun(VT200_0851_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nwait(2);`nTakeScreenshot(VT200-0851);`nvalidate4;"
$ws.Range("H2").Value = "validate1`n{`nvalidate_PageTitle=Compliance JS specs`n};`nvalidate2`n{`nvalidate_PageTitle=Signal JS Test`n};`nvalidate3`n{`nvalidate_Text_Exists=VT200-0851`n};`nvalidate4`n{`nvalidate_Screenshot=VT200-0851`nvalidate_Iconposition=signalview_xpath,left,20`nvalidate_Iconposition=signalview_xpath,top,40`n};"

$ws.Range("G5").Value = "wait(3);`nvalidate1;`nlink_Click(signal_test_link);`nvalidate2;`nSelectTestToRun(VT200_0854_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nwait(2);`nvalidate4;"
$ws.Range("H5").Value = "validate1`n{`nvalidate_PageTitle=Compliance JS specs`n};`nvalidate2`n{`nvalidate_PageTitle=Signal JS Test`n};`nvalidate3`n{`nvalidate_Text_Exists=VT200-0854`n};`nvalidate4`n{`nvalidate_Iconposition=signalview_xpath,left,40`n};"

$ws.Range("G6").Value = "wait(3);`nvalidate1;`nlink_Click(signal_test_link);`nvalidate2;`nSelectTestToRun(VT200_0855_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nwait(2);`nvalidate4;"
$ws.Range("H6").Value = "validate1`n{`nvalidate_PageTitle=Compliance JS specs`n};`nvalidate2`n{`nvalidate_PageTitle=Signal JS Test`n};`nvalidate3`n{`nvalidate_Text_Exists=VT200-0855`n};`nvalidate4`n{`nvalidate_Iconposition=signalview_xpath,top,40`n};"

$ws.Range("G12").Value = "wait(3);`nvalidate1;`nlink_Click(signal_test_link);`nvalidate2;`nSelectTestToRun(VT200_0861_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nwait(2);`nvalidate4;`nwait(12);`nvalidate5;"
$ws.Range("H12").Value = "validate1`n{`nvalidate_PageTitle=Compliance JS specs`n};`nvalidate2`n{`nvalidate_PageTitle=Signal JS Test`n};`nvalidate3`n{`nvalidate_Text_Exists=VT200-0861`n};`nvalidate4`n{`nvalidate_isIconDisplayed=signalview_xpath,true`n};`nvalidate5`n{`nvalidate_isIconDisplayed=signalview_xpath,false`n};"

$ws.Rows.Item(2).RowHeight = 203.25

$ws.Range("G2").Select()

